$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. List item "C#" / " " -> append a new run "/ .Net Core 3.1" right after it
#    (keeps the existing "C#" and " " runs untouched, inserts a brand new run).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("C# ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Collapse(0)
$r1.InsertAfter("/ .Net Core 3.1")
$r1.LanguageID = "en-US"
$r1.Bold = 1
$r1.Bold = 0

# ---------------------------------------------------------------------------
# 2. Merge "Do some configuration changes ... e.g. <" into a single run,
#    dropping the proofErr grammar markup that used to wrap "e.g.".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Do some configuration changes to make the project compatible for Azure DevOps and chrome e.g. <",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Do some configuration changes to make the project compatible for Azure DevOps and chrome e.g. <", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Merge "file can be integrated in pipeline so people have used their own
#    different ways to do it." (drop proofErr around "pipeline").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " file can be integrated in pipeline so people have used their own different ways to do it.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " file can be integrated in pipeline so people have used their own different ways to do it.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Merge "   One of the way is given below.  " (drop proofErr around "way").
#    This sits right after the run from step 3 with identical formatting, so
#    nudge it (harmless Bold toggle) to keep it a distinct run instead of
#    silently folding back into the previous one.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute(
    "   One of the way is given below.  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "   One of the way is given below.  ", 2) | Out-Null
$r4.Bold = 1
$r4.Bold = 0

# ---------------------------------------------------------------------------
# 5. Merge "One of the way is: " (drop proofErr around "way").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "One of the way is: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "One of the way is: ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6. Merge " files i.e. run.sh that launches Docker image and run the tests
#    and create test.sh for each " (drop proofErr around "i.e.").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " files i.e. run.sh that launches Docker image and run the tests and create test.sh for each ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " files i.e. run.sh that launches Docker image and run the tests and create test.sh for each ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7. "Please note that ... was just explain the process ... definitely able
#    to do it and integrate " paragraph:
#      a) merge the tail ("explain the process ... definitely able to do it
#         and integrate ") into a single run, dropping proofErr around
#         "definitely";
#      b) insert a brand-new run "to " right after "...was just " (forcing a
#         distinct run via a harmless Bold toggle so it doesn't get folded
#         back into its identically-formatted neighbours).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "explain the process but if I have to do it in actual, I can definitely able to do it and integrate ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "explain the process but if I have to do it in actual, I can definitely able to do it and integrate ", 2) | Out-Null

$r7 = $d.Content
$r7.Find.Execute("was just ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $d.Range($r7.End, $r7.End)
$ins.InsertAfter("to ")
$ins.LanguageID = "en-US"
$ins.Bold = 1
$ins.Bold = 0
